# "add default config for npc"
# Rename the "Player" NPC config row to "Player_0_0", and add a new
# "Default" NPC config row (row 7) cloned from the existing melee-NPC
# stat block (same stats as AttackNpc3 / Player_0_0 / Enemy rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Player -> Player_0_0
$ws.Range("A5").Value = "Player_0_0"

# Row 7: new "Default" NPC entry
$ws.Range("A7").Value = "Default"
$ws.Range("B7").Value = "60"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = "0"
$ws.Range("E7").Value = "60"
$ws.Range("F7").Value = 600
$ws.Range("G7").Value = 600
$ws.Range("H7").Value = 600
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 150
$ws.Range("L7").Value = 80
$ws.Range("M7").Value = 55000
$ws.Range("N7").Value = 10000
$ws.Range("O7").Value = 10
$ws.Range("P7").Value = 10
$ws.Range("Q7").Value = 10
$ws.Range("R7").Value = 10
$ws.Range("S7").Value = 5
$ws.Range("T7").Value = 5
$ws.Range("U7").Value = 5
$ws.Range("V7").Value = 5
$ws.Range("W7").Value = "Prefabs/Object/Alch_plate"
$ws.Range("X7").Value = 2
$ws.Range("Y7").Value = 20
$ws.Range("Z7").Value = "DropBag_1"

# Widen the Prefab column (W) to fit its content.
$ws.Columns.Item(23).ColumnWidth = 62.7

# Restore the cursor/selection to A6 (also clears the stale topLeftCell
# scroll position left over from the previous save).
[void]$ws.Range("A6").Select()

Write-Output "done"
